$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.03086786270141601, 0.05265259742736816, 0.01334075927734375, 0.03766536712646484, 0.004376602172851562, 0.1255281448364258, 0.02653508186340332, 0.0518836498260498, 0.01384758949279785, 0.02770719528198242, 0.007568073272705078, 0.02505688667297363)
    3  = @(0.1376707077026367, 0.04104599952697754, 0.02620229721069336, 0.005586051940917968, 0.01305994987487793, 0.01526103019714356, 0.1200681686401367, 0.03180246353149414, 0.0590017318725586, 0.02356057167053223, 0.03096041679382324, 0.01335206031799316)
    4  = @(0.05341620445251465, 0.02685074806213379, 0.009375047683715821, 0.006252241134643555, 0.04615321159362793, 0.00760345458984375, 0.02696304321289062, 0.02380952835083008, 0.01371746063232422, 0.01717348098754883, 0.0430755615234375, 0.008268213272094727)
    5  = @(0.02996811866760254, 0.02252864837646484, 0.011602783203125, 0.0157376766204834, 0, 0, 0.01835823059082031, 0.02276477813720703, 0.01563777923583985, 0.02546634674072266, 0, 0)
    6  = @(0.4853501319885254, 0.03035573959350586, 0.4424353122711182, 0.02622151374816895, 0.4698287010192871, 0.01604518890380859, 0.1610980987548828, 0.01732673645019531, 0.1611368656158447, 0.01783418655395508, 0.5247678279876709, 0.01622662544250488)
    7  = @(1.041378164291382, 0.1161933422088623, 0.3512662887573242, 0.03590364456176758, 0.7601140975952149, 0.04442977905273438, 0.9165266513824463, 0.09882216453552246, 0.4718425273895264, 0.06103196144104004, 0.7493313789367676, 0.03231759071350097)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
